# Devon Air Ambulance workbook restructure:
#  - rename "landing" -> "welcome" and add welcome/coming-soon text
#  - insert a brand new "what_is" sheet right after "welcome"
#  - leave "setup" / "scenario" / "output" sheets untouched (they simply
#    shift right in tab order, which Excel/the exporter handles for us)

$wb = $excel.ActiveWorkbook

# --- 1. "landing" -> "welcome" -------------------------------------------------
$welcome = $wb.Worksheets.Item("landing")
$welcome.Name = "welcome"

# Fill in the new text cells (order matches how the shared-string table
# ends up populated: welcome sheet first, then the new sheet's A column,
# then the new sheet's B column).
$welcome.Range("B2").Value = "Welcome to the Devon Air Ambulance Simulation Model"
$welcome.Range("B3").Value = "Coming Soon!"

# --- 2. brand new "what_is" sheet, placed right after "welcome" ---------------
$whatIs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $welcome)
$whatIs.Name = "what_is"

$whatIs.Range("A1").Value = "reference"
$whatIs.Range("A2").Value = "page_title"
$whatIs.Range("A3").Value = "page_description"
$whatIs.Range("A4").Value = "tab_1_name"
$whatIs.Range("A5").Value = "tab_2_name"
$whatIs.Range("A6").Value = "tab_3_name"
$whatIs.Range("A7").Value = "tab_1_content"
$whatIs.Range("A8").Value = "tab_2_content"
$whatIs.Range("A9").Value = "tab_3_content"

$whatIs.Range("B1").Value = "text"
$whatIs.Range("B2").Value = "An Introduction to Discrete Event Simulation Modelling"
$whatIs.Range("B3").Value = "This page will introduce you to the fundamental concepts of the computer simulation techniques used in this model.`nNo prior knowledge of computer simulation, maths or data science will be required."
$whatIs.Range("B4").Value = "An Introduction to Simulation Modelling"
$whatIs.Range("B5").Value = "Benefits of Simulation Modelling"
$whatIs.Range("B6").Value = "Limitations of Simulation Modelling"
$whatIs.Range("B7").Value = "Coming Soon!"
$whatIs.Range("B8").Value = "Coming Soon!"
$whatIs.Range("B9").Value = "Coming Soon!"

# Formatting to match the "what_is" layout: wide text column, wrapped
# description cell with a taller row.
$whatIs.Columns.Item(1).ColumnWidth = 17.92
$whatIs.Columns.Item(2).ColumnWidth = 70.26
$whatIs.Range("B3").WrapText = $true
$whatIs.Rows.Item(3).RowHeight = 52.8

# --- 3. selection / active-sheet state -----------------------------------------
$welcome.Activate()
$welcome.Range("A1:B1").Select() | Out-Null

$whatIs.Activate()
$whatIs.Range("C7").Select() | Out-Null
